$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"5086.153335856128"
$ws.Range("H2").Value = [double]"210"
$ws.Range("I2").Value = [double]"786"
$ws.Range("J2").Value = [double]"996"
$ws.Range("K2").Value = [double]"0.7891566265060241"

$ws.Range("E3").Value = [double]"5086.153335856128"
$ws.Range("H3").Value = [double]"11580"
$ws.Range("I3").Value = [double]"666"
$ws.Range("J3").Value = [double]"12246"
$ws.Range("K3").Value = [double]"0.05438510534051935"

$ws.Range("E4").Value = [double]"1204.914932712188"
$ws.Range("F4").Value = [double]"1.439258204408449e-247"
$ws.Range("H4").Value = [double]"257"
$ws.Range("J4").Value = [double]"476"
$ws.Range("K4").Value = [double]"0.4600840336134454"

$ws.Range("E5").Value = [double]"1204.914932712188"
$ws.Range("F5").Value = [double]"1.439258204408449e-247"
$ws.Range("H5").Value = [double]"66"
$ws.Range("J5").Value = [double]"100"
$ws.Range("K5").Value = [double]"0.34"

$ws.Range("E6").Value = [double]"1204.914932712188"
$ws.Range("F6").Value = [double]"1.439258204408449e-247"
$ws.Range("I6").Value = [double]"57"
$ws.Range("J6").Value = [double]"186"
$ws.Range("K6").Value = [double]"0.3064516129032258"

$ws.Range("E7").Value = [double]"1204.914932712188"
$ws.Range("F7").Value = [double]"1.439258204408449e-247"
$ws.Range("I7").Value = [double]"50"
$ws.Range("J7").Value = [double]"187"
$ws.Range("K7").Value = [double]"0.267379679144385"

$ws.Range("E8").Value = [double]"1204.914932712188"
$ws.Range("F8").Value = [double]"1.439258204408449e-247"
$ws.Range("H8").Value = [double]"118"
$ws.Range("J8").Value = [double]"156"
$ws.Range("K8").Value = [double]"0.2435897435897436"

$ws.Range("E9").Value = [double]"1204.914932712188"
$ws.Range("F9").Value = [double]"1.439258204408449e-247"
$ws.Range("H9").Value = [double]"107"
$ws.Range("J9").Value = [double]"141"
$ws.Range("K9").Value = [double]"0.2411347517730496"

$ws.Range("E10").Value = [double]"1204.914932712188"
$ws.Range("F10").Value = [double]"1.439258204408449e-247"

$ws.Range("E11").Value = [double]"1204.914932712188"
$ws.Range("F11").Value = [double]"1.439258204408449e-247"
$ws.Range("H11").Value = [double]"146"
$ws.Range("J11").Value = [double]"186"
$ws.Range("K11").Value = [double]"0.2150537634408602"

$ws.Range("E12").Value = [double]"1204.914932712188"
$ws.Range("F12").Value = [double]"1.439258204408449e-247"
$ws.Range("H12").Value = [double]"599"
$ws.Range("I12").Value = [double]"131"
$ws.Range("K12").Value = [double]"0.1794520547945206"

$ws.Range("E13").Value = [double]"1204.914932712188"
$ws.Range("F13").Value = [double]"1.439258204408449e-247"
$ws.Range("H13").Value = [double]"834"
$ws.Range("I13").Value = [double]"149"
$ws.Range("J13").Value = [double]"983"
$ws.Range("K13").Value = [double]"0.1515768056968464"

$ws.Range("E14").Value = [double]"1204.914932712188"
$ws.Range("F14").Value = [double]"1.439258204408449e-247"
$ws.Range("H14").Value = [double]"833"
$ws.Range("J14").Value = [double]"941"
$ws.Range("K14").Value = [double]"0.1147715196599362"

$ws.Range("E15").Value = [double]"1204.914932712188"
$ws.Range("F15").Value = [double]"1.439258204408449e-247"
$ws.Range("H15").Value = [double]"1026"
$ws.Range("I15").Value = [double]"109"
$ws.Range("J15").Value = [double]"1135"
$ws.Range("K15").Value = [double]"0.0960352422907489"

$ws.Range("E16").Value = [double]"1204.914932712188"
$ws.Range("F16").Value = [double]"1.439258204408449e-247"
$ws.Range("H16").Value = [double]"1145"
$ws.Range("I16").Value = [double]"95"
$ws.Range("K16").Value = [double]"0.07661290322580645"

$ws.Range("E17").Value = [double]"1204.914932712188"
$ws.Range("F17").Value = [double]"1.439258204408449e-247"
$ws.Range("H17").Value = [double]"1482"
$ws.Range("J17").Value = [double]"1579"
$ws.Range("K17").Value = [double]"0.06143128562381254"

$ws.Range("E18").Value = [double]"1204.914932712188"
$ws.Range("F18").Value = [double]"1.439258204408449e-247"
$ws.Range("H18").Value = [double]"2517"
$ws.Range("I18").Value = [double]"123"
$ws.Range("J18").Value = [double]"2640"
$ws.Range("K18").Value = [double]"0.04659090909090909"

$ws.Range("E19").Value = [double]"1204.914932712188"
$ws.Range("F19").Value = [double]"1.439258204408449e-247"
$ws.Range("H19").Value = [double]"2153"
$ws.Range("J19").Value = [double]"2252"
$ws.Range("K19").Value = [double]"0.04396092362344583"

$ws.Range("E20").Value = [double]"915.4726472302203"
$ws.Range("F20").Value = [double]"4.248980687842385e-201"
$ws.Range("H20").Value = [double]"1198"
$ws.Range("I20").Value = [double]"563"
$ws.Range("J20").Value = [double]"1761"
$ws.Range("K20").Value = [double]"0.3197047132311187"

$ws.Range("E21").Value = [double]"915.4726472302203"
$ws.Range("F21").Value = [double]"4.248980687842385e-201"
$ws.Range("H21").Value = [double]"10592"
$ws.Range("I21").Value = [double]"889"
$ws.Range("J21").Value = [double]"11481"
$ws.Range("K21").Value = [double]"0.07743227941816914"

$ws.Range("E22").Value = [double]"392.3837685207457"
$ws.Range("F22").Value = [double]"2.752934587762642e-74"
$ws.Range("H22").Value = [double]"128"
$ws.Range("I22").Value = [double]"85"
$ws.Range("J22").Value = [double]"213"
$ws.Range("K22").Value = [double]"0.3990610328638498"

$ws.Range("E23").Value = [double]"392.3837685207457"
$ws.Range("F23").Value = [double]"2.752934587762642e-74"
$ws.Range("H23").Value = [double]"210"
$ws.Range("J23").Value = [double]"290"
$ws.Range("K23").Value = [double]"0.2758620689655172"

$ws.Range("E24").Value = [double]"392.3837685207457"
$ws.Range("F24").Value = [double]"2.752934587762642e-74"
$ws.Range("I24").Value = [double]"165"
$ws.Range("J24").Value = [double]"909"
$ws.Range("K24").Value = [double]"0.1815181518151815"

$ws.Range("E25").Value = [double]"392.3837685207457"
$ws.Range("F25").Value = [double]"2.752934587762642e-74"
$ws.Range("H25").Value = [double]"4400"
$ws.Range("I25").Value = [double]"516"
$ws.Range("J25").Value = [double]"4916"
$ws.Range("K25").Value = [double]"0.1049633848657445"

$ws.Range("E26").Value = [double]"392.3837685207457"
$ws.Range("F26").Value = [double]"2.752934587762642e-74"
$ws.Range("H26").Value = [double]"6308"
$ws.Range("I26").Value = [double]"606"
$ws.Range("J26").Value = [double]"6914"
$ws.Range("K26").Value = [double]"0.08764824992768296"

$ws.Range("E27").Value = [double]"335.8134709170386"
$ws.Range("F27").Value = [double]"6.459790618508831e-67"
$ws.Range("H27").Value = [double]"3088"
$ws.Range("I27").Value = [double]"651"
$ws.Range("J27").Value = [double]"3739"
$ws.Range("K27").Value = [double]"0.1741107247927253"

$ws.Range("E28").Value = [double]"335.8134709170386"
$ws.Range("F28").Value = [double]"6.459790618508831e-67"
$ws.Range("H28").Value = [double]"4732"
$ws.Range("I28").Value = [double]"602"
$ws.Range("J28").Value = [double]"5334"
$ws.Range("K28").Value = [double]"0.1128608923884514"

$ws.Range("E29").Value = [double]"335.8134709170386"
$ws.Range("F29").Value = [double]"6.459790618508831e-67"
$ws.Range("H29").Value = [double]"154"
$ws.Range("J29").Value = [double]"164"
$ws.Range("K29").Value = [double]"0.06097560975609756"

$ws.Range("E30").Value = [double]"335.8134709170386"
$ws.Range("F30").Value = [double]"6.459790618508831e-67"
$ws.Range("H30").Value = [double]"716"
$ws.Range("J30").Value = [double]"760"
$ws.Range("K30").Value = [double]"0.05789473684210526"

$ws.Range("E31").Value = [double]"335.8134709170386"
$ws.Range("F31").Value = [double]"6.459790618508831e-67"
$ws.Range("H31").Value = [double]"400"
$ws.Range("J31").Value = [double]"421"
$ws.Range("K31").Value = [double]"0.0498812351543943"

$ws.Range("E32").Value = [double]"335.8134709170386"
$ws.Range("F32").Value = [double]"6.459790618508831e-67"
$ws.Range("H32").Value = [double]"2265"
$ws.Range("J32").Value = [double]"2376"
$ws.Range("K32").Value = [double]"0.04671717171717172"

$ws.Range("E33").Value = [double]"335.8134709170386"
$ws.Range("F33").Value = [double]"6.459790618508831e-67"
$ws.Range("H33").Value = [double]"435"
$ws.Range("J33").Value = [double]"448"
$ws.Range("K33").Value = [double]"0.02901785714285714"

$ws.Range("E34").Value = [double]"109.391853707005"
$ws.Range("F34").Value = [double]"1.331757529725758e-25"
$ws.Range("H34").Value = [double]"6864"
$ws.Range("I34").Value = [double]"1053"
$ws.Range("J34").Value = [double]"7917"
$ws.Range("K34").Value = [double]"0.1330049261083744"

$ws.Range("E35").Value = [double]"109.391853707005"
$ws.Range("F35").Value = [double]"1.331757529725758e-25"
$ws.Range("H35").Value = [double]"4926"
$ws.Range("I35").Value = [double]"399"
$ws.Range("J35").Value = [double]"5325"
$ws.Range("K35").Value = [double]"0.07492957746478873"

$ws.Range("E36").Value = [double]"86.00805738538253"
$ws.Range("F36").Value = [double]"1.792039823210119e-20"
$ws.Range("H36").Value = [double]"7714"
$ws.Range("I36").Value = [double]"1127"
$ws.Range("J36").Value = [double]"8841"
$ws.Range("K36").Value = [double]"0.1274742676167854"

$ws.Range("E37").Value = [double]"86.00805738538253"
$ws.Range("F37").Value = [double]"1.792039823210119e-20"
$ws.Range("H37").Value = [double]"4076"
$ws.Range("I37").Value = [double]"325"
$ws.Range("J37").Value = [double]"4401"
$ws.Range("K37").Value = [double]"0.07384685298795728"

$ws.Range("E38").Value = [double]"42.314501917315"
$ws.Range("F38").Value = [double]"7.771416267429973e-11"
$ws.Range("H38").Value = [double]"1340"
$ws.Range("I38").Value = [double]"251"
$ws.Range("K38").Value = [double]"0.1577624135763671"

$ws.Range("E39").Value = [double]"42.314501917315"
$ws.Range("F39").Value = [double]"7.771416267429973e-11"
$ws.Range("H39").Value = [double]"10450"
$ws.Range("I39").Value = [double]"1201"
$ws.Range("J39").Value = [double]"11651"
$ws.Range("K39").Value = [double]"0.1030812805767745"

$ws.Range("E40").Value = [double]"34.81280277953115"
$ws.Range("F40").Value = [double]"3.629786486834646e-09"
$ws.Range("H40").Value = [double]"6020"
$ws.Range("I40").Value = [double]"861"
$ws.Range("J40").Value = [double]"6881"
$ws.Range("K40").Value = [double]"0.1251271617497457"

$ws.Range("E41").Value = [double]"34.81280277953115"
$ws.Range("F41").Value = [double]"3.629786486834646e-09"
$ws.Range("H41").Value = [double]"5770"
$ws.Range("I41").Value = [double]"591"
$ws.Range("J41").Value = [double]"6361"
$ws.Range("K41").Value = [double]"0.09290991982392706"
